$p = $ppt.ActivePresentation
$EMU = 12700.0

# --- Slide 4 (Course Model - Side View) ---
$s4 = $p.Slides.Item(4)
$shape = $s4.Shapes.AddTextbox(1, 1859863/$EMU, 3148346/$EMU, 914417/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Magnet"

$shape = $s4.Shapes.AddTextbox(1, 3324225/$EMU, 3059389/$EMU, 1165191/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Electrodes"

$shape = $s4.Shapes.AddTextbox(1, 4407214/$EMU, 5681212/$EMU, 1165191/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Electrodes"

$shape = $s4.Shapes.AddTextbox(1, 329346/$EMU, 3676649/$EMU, 457176/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Lid"

$shape = $s4.Shapes.AddTextbox(1, 330864/$EMU, 5310971/$EMU, 684803/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Basin"

$shape = $s4.Shapes.AddTextbox(1, 5380179/$EMU, 4540315/$EMU, 1063304/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Transwell"

$shape = $s4.Shapes.AddTextbox(1, 1661043/$EMU, 5190836/$EMU, 914417/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Magnet"

# --- Slide 5 (Course Model - Top View) ---
$s5 = $p.Slides.Item(5)
$shape = $s5.Shapes.AddTextbox(1, 4545584/$EMU, 3657560/$EMU, 551754/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "PCB"

$shape = $s5.Shapes.AddTextbox(1, 4697984/$EMU, 3809960/$EMU, 551754/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "PCB"

$shape = $s5.Shapes.AddTextbox(1, 9128582/$EMU, 3640781/$EMU, 1205908/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Receptacle"

$shape = $s5.Shapes.AddTextbox(1, 2541616/$EMU, 2755857/$EMU, 725070/$EMU, 369332/$EMU)
$shape.TextFrame.WordWrap = $false
$shape.TextFrame.AutoSize = 1
$shape.Fill.Visible = $false
$shape.TextFrame.TextRange.Text = "Wires"
